# Populate the "Membre du groupe" table (Prénom / Nom) on Sheet1.
# Cell values are written in the same order the original author typed
# them in Excel, so the generated sharedStrings table lines up with the
# source file (header row, then first-name column of row 2, then the
# rest row-by-row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Prénom"
$ws.Range("B1").Value = "Nom"

$ws.Range("B2").Value = "KHADIDIATOU"
$ws.Range("A2").Value = "COULIBALY"

$ws.Range("A3").Value = "NDONG"
$ws.Range("B3").Value = "TAMSIR"

$ws.Range("B4").Value = "SAMBA"
$ws.Range("A4").Value = "DIENG"

$ws.Range("A5").Value = "ONANENA AMANA"
$ws.Range("B5").Value = "JEANNE DE LA FLECHE"

# Auto-fit the "Prénom" column to the longest entry, like the author did
# after typing the data in.
$ws.Columns.Item(1).AutoFit()

# Leave the selection where the author ended up (cell C8) after entering
# the table.
$ws.Range("C8").Select() | Out-Null
